# Update the cryptos list sheet with refreshed Price / Volume(1h) figures,
# matching the GitHub Actions scheduled data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force the cell to be treated as text (the sheet stores prices/volumes
    # as plain text, e.g. "27.434.55" or "1.002") rather than letting Excel
    # auto-convert numeric-looking strings into real numbers. Resetting the
    # style back to "Normal" afterwards avoids leaving a stray number-format
    # style behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Rows where only the Price (D) and Volume(1h) (E) columns changed.
$simpleUpdates = @(
    @{ Row = 2;  D = "27.434.55";    E = "  -2.46%  " },
    @{ Row = 3;  D = "1.746.17";     E = "  -2.82%  " },
    @{ Row = 4;  D = "1.002";        E = "  -0.31%  " },
    @{ Row = 5;  D = "322.78";       E = "  -3.48%  " },
    @{ Row = 6;  D = "0.9995";       E = "  -0.11%  " },
    @{ Row = 7;  D = "0.4254";       E = "  -8.78%  " },
    @{ Row = 8;  D = "0.3627";       E = "  -2.34%  " },
    @{ Row = 9;  D = "45.19";        E = "  +0.48%  " },
    @{ Row = 10; D = "0.07441";      E = "  -2.13%  " },
    @{ Row = 11; D = "1.118";        E = "  -2.64%  " },
    @{ Row = 12; D = "0.9999";       E = "  -0.40%  " },
    @{ Row = 13; D = "21.57";        E = "  -2.98%  " },
    @{ Row = 14; D = "6.104";        E = "  -3.58%  " },
    @{ Row = 15; D = "7.201";        E = "  -2.22%  " },
    @{ Row = 16; D = "1.743.48";     E = "  -2.91%  " },
    @{ Row = 17; D = "0.00001070";   E = "  -1.60%  " },
    @{ Row = 18; D = "88.02";        E = "  +7.13%  " },
    @{ Row = 19; D = "0.06022";      E = "  -10.28%  " },
    @{ Row = 20; D = "0.9994";       E = "  -0.13%  " },
    @{ Row = 21; D = "16.90";        E = "  -2.66%  " },
    @{ Row = 22; D = "6.114";        E = "  -4.18%  " },
    @{ Row = 23; D = "0.5238";       E = "  -5.23%  " },
    @{ Row = 24; D = "27.464.37";    E = "  -2.41%  " },
    @{ Row = 25; D = "11.40";        E = "  -3.98%  " },
    @{ Row = 26; D = "2.376";        E = "  -1.38%  " },
    @{ Row = 29; D = "150.62";       E = "  -0.68%  " },
    @{ Row = 30; D = "1.939.15";     E = "  -3.46%  " },
    @{ Row = 31; D = "126.59";       E = "  -5.10%  " },
    @{ Row = 32; D = "1.192";        E = "  -4.10%  " },
    @{ Row = 33; D = "5.725";        E = "  -2.18%  " },
    @{ Row = 34; D = "0.09112";      E = "  -5.33%  " },
    @{ Row = 35; D = "3.585";        E = "  -11.28%  " },
    @{ Row = 36; D = "12.99";        E = "  +7.37%  " },
    @{ Row = 37; D = "0.2148";       E = "  -3.22%  " },
    @{ Row = 38; D = "5.094";        E = "  -2.46%  " },
    @{ Row = 39; D = "0.02269";      E = "  -3.90%  " },
    @{ Row = 42; D = "1.189";        E = "  -3.17%  " },
    @{ Row = 43; D = "8.063";        E = "  -0.18%  " },
    @{ Row = 46; D = "13.67";        E = "  -3.44%  " },
    @{ Row = 47; D = "3.720";        E = "  -2.64%  " },
    @{ Row = 48; D = "0.5860";       E = "  -4.00%  " },
    @{ Row = 49; D = "125.16";       E = "  -3.08%  " },
    @{ Row = 50; D = "1.962";        E = "  -4.22%  " },
    @{ Row = 51; D = "0.06866";      E = "  -3.93%  " }
)

foreach ($u in $simpleUpdates) {
    Set-TextCell $ws.Cells.Item($u.Row, 4) $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Rows where the ranking reshuffled, so the coin (B), its link (C), and its
# own Price/Volume (D/E) all moved to a different row.
$swapUpdates = @(
    @{ Row = 27; B = "LidoDAOToken";    C = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo";             D = "2.409";   E = "  +1.75%  " },
    @{ Row = 28; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc";      D = "20.36";   E = "  -1.47%  " },
    @{ Row = 40; B = "TheSandbox";      C = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand";              D = "0.6434";  E = "  -3.40%  " },
    @{ Row = 41; B = "Hedera";          C = "https://coinranking.com/coin/jad286TjB+hedera-hbar";                  D = "0.06068"; E = "  -4.41%  " },
    @{ Row = 44; B = "WEMIXTOKEN";      C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix";             D = "1.434";   E = "  -6.62%  " },
    @{ Row = 45; B = "Frax";            C = "https://coinranking.com/coin/KfWtaeV1W+frax-frax";                    D = "0.9987";  E = "  -0.06%  " }
)

foreach ($u in $swapUpdates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    Set-TextCell $ws.Cells.Item($u.Row, 4) $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
